# Commit: Fri, Jul 24, 2020 12:05:17 AM
#
# 1) Swap the presentation's applied theme colour scheme from the
#    "Integral" / "Red Violet" palette over to the stock "Office Theme" /
#    "Office" palette (this is what happens in the OOXML when the Design
#    gallery selection is changed back to the default Office theme).
# 2) Re-point the three data tables (slides 14, 15 and 16) from the
#    deck's custom "Table_0" table style onto the built-in table style
#    {BCA6F262-E643-46AB-AC83-BF9520C70F70}.

$p = $ppt.ActivePresentation

# --- 1. Theme colours -------------------------------------------------
# Target palette = the "Office Theme" / "Office" colour scheme
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), expressed as the
# BGR-packed COM RGB() values PowerPoint's ThemeColorScheme expects.
$m = $p.SlideMaster
$colors = $m.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0x000000   # dk1      000000
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colors.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colors.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colors.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colors.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$colors.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$colors.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$colors.Colors(11).RGB = 0xC16305   # hlink    0563C1
$colors.Colors(12).RGB = 0x724F95   # folHlink 954F72

# --- 2. Table styles ----------------------------------------------------
$targetStyle = "{BCA6F262-E643-46AB-AC83-BF9520C70F70}"
foreach ($idx in 14, 15, 16) {
    $s  = $p.Slides.Item($idx)
    $sh = $s.Shapes.Item(1)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle($targetStyle)
    }
}
